$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Excel file total marks error": update the marking scheme and totals
# for row 11 (Marking) and row 12 (Total) on the marksheet.

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 108
$ws.Range("E12").Value = "106 / 112"
